$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for new columns I and J, matching style of existing headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in I0 and IF values for each data row (2-45)
$data = @(
    @(7,8),
    @(8,8),
    @(7,8),
    @(6,7),
    @(6,7),
    @(6,7),
    @(6,7),
    @(8,8),
    @(7,7),
    @(6,6),
    @(8,8),
    @(5,5),
    @(6,6),
    @(8,8),
    @(7,7),
    @(7,8),
    @(7,7),
    @(6,7),
    @(7,7),
    @(6,6),
    @(6,7),
    @(6,6),
    @(6,7),
    @(6,6),
    @(7,7),
    @(6,6),
    @(8,8),
    @(5,5),
    @(7,7),
    @(6,6),
    @(1,2),
    @(7,8),
    @(6,7),
    @(6,6),
    @(7,8),
    @(4,5),
    @(7,8),
    @(9,9),
    @(6,8),
    @(5,5),
    @(8,8),
    @(9,9),
    @(4,5),
    @(9,9)
)

for ($k = 0; $k -lt $data.Length; $k++) {
    $row = 2 + $k
    $ws.Cells.Item($row, 9).Value = $data[$k][0]
    $ws.Cells.Item($row, 10).Value = $data[$k][1]
}

